# Fixed update to excel issue
# - Rename header columns to clarify PO quantity metrics
# - Add a new "PO Forecast" sheet with the Prophet-style forecast output

$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last (3rd) sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header / date formatting from the Weekly Quantity sheet so the
# new sheet matches the look of the existing ones (bold header w/ border,
# centered + top aligned, and the yyyy-mm-dd date style on column A).
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A41").PasteSpecial(-4122)

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows 2-41
$poForecastData = @(
    @(45298.99999999999, 76, 25.48247646433178, 126.8910763593615),
    @(45305.99999999999, 75, 23.18912203980951, 126.6978431833387),
    @(45312.99999999999, 74, 17.51843648222154, 122.8670395250111),
    @(45326.99999999999, 71, 19.72771582953012, 125.1754890777817),
    @(45340.99999999999, 68, 15.91141812573505, 119.4016146033349),
    @(45347.99999999999, 67, 17.60622590519688, 119.0117822357901),
    @(45354.99999999999, 66, 10.83091834862257, 118.4258978608859),
    @(45368.99999999999, 63, 7.219632921571003, 114.1253776110131),
    @(45375.99999999999, 62, 9.664660441445772, 115.0661050518691),
    @(45382.99999999999, 60, 10.19266483395189, 111.5751251402444),
    @(45389.99999999999, 59, 10.61028087013815, 110.5781168488856),
    @(45417.99999999999, 54, 4.231860368505822, 104.2429130478735),
    @(45424.99999999999, 53, 0.9014568875121129, 107.0611034834778),
    @(45431.99999999999, 51, -3.346076095389446, 104.5292568076603),
    @(45438.99999999999, 50, -3.848859019852029, 100.7900127513933),
    @(45452.99999999999, 47, -4.731608790128378, 104.0122278782882),
    @(45466.99999999999, 45, -13.23221591873498, 98.46887136755896),
    @(45473.99999999999, 43, -12.0837724398032, 96.01243637368039),
    @(45494.99999999999, 39, -15.81026262742467, 90.54824233061497),
    @(45508.99999999999, 37, -17.68967060081302, 89.46788692744225),
    @(45515.99999999999, 35, -18.10630374280814, 90.47804300436731),
    @(45522.99999999999, 34, -12.80665731079499, 83.37566376680283),
    @(45529.99999999999, 33, -20.27762308029343, 83.16093221108983),
    @(45536.99999999999, 31, -21.13314314366059, 81.70804569393398),
    @(45543.99999999999, 30, -20.70072558778118, 82.28130341594672),
    @(45550.99999999999, 29, -24.13200144414039, 78.28125518741582),
    @(45557.99999999999, 27, -25.66326870469888, 84.56913201923048),
    @(45564.99999999999, 26, -29.31862581591865, 79.42200651380045),
    @(45571.99999999999, 25, -27.67098447716745, 77.19374414690022),
    @(45585.99999999999, 22, -29.23034602616552, 71.37090640334918),
    @(45592.99999999999, 21, -31.68455451743926, 73.9201533832731),
    @(45599.99999999999, 20, -29.32708525012992, 73.19003326138481),
    @(45606.99999999999, 18, -37.08152257014945, 70.5026735473494),
    @(45613.99999999999, 17, -37.82135944421326, 68.41352151550993),
    @(45620.99999999999, 16, -36.5094444720025, 70.138109104474),
    @(45627.99999999999, 14, -39.75675596213006, 67.04333447392159),
    @(45634.99999999999, 13, -37.9857576371945, 65.86238130665437),
    @(45641.99999999999, 12, -40.67424950694421, 65.15285175528565),
    @(45648.99999999999, 10, -40.01072762862336, 63.86812182287215),
    @(45655.99999999999, 9, -46.1692383626399, 55.67068384812487)
)

$r = 2
foreach ($row in $poForecastData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

Write-Output "PO Forecast sheet added; headers updated."
